{"js": "// The commit removes the italic Hindi subtitle paragraph (\"\u0930\u0942\u0924\") that\n// immediately follows the \"RUT\" Heading 2 paragraph at the top of the\n// book-intro section, merging what is left into the \"RUT\" heading\n// paragraph (the trailing empty run/paragraph mark of the subtitle\n// paragraph is what is kept).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the unique \"RUT\" Heading 2 paragraph (the book-code heading).\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Heading 2\" && para.text === \"RUT\") {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex === -1) {\n  throw new Error('Could not find the \"RUT\" Heading 2 paragraph.');\n}\n\n// The very next paragraph is the italic Hindi subtitle (\"\u0930\u0942\u0924\") that the\n// commit deletes (along with the paragraph break that separated it from\n// the heading above).\nconst subtitlePara = paragraphs.items[headingIndex + 1];\nsubtitlePara.delete();\nawait context.sync();\n", "ps1": "# The commit removes the italic Hindi subtitle paragraph (\"\u0930\u0942\u0924\") that\n# immediately follows the \"RUT\" Heading 2 paragraph at the top of the\n# book-intro section, merging what is left into the \"RUT\" heading\n# paragraph (the trailing empty run/paragraph mark of the subtitle\n# paragraph is what is kept).\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $trimmed = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($trimmed -eq \"RUT\" -and $p.Style.NameLocal -eq \"Heading 2\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the 'RUT' Heading 2 paragraph.\"\n}\n\n# The very next paragraph is the italic Hindi subtitle (\"\u0930\u0942\u0924\") that the\n# commit deletes (along with the paragraph break that separated it from\n# the heading above).\n$subtitlePara = $d.Paragraphs.Item($targetIndex + 1)\n$subtitlePara.Range.Delete()\n"}
